$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Pandey"
$ws.Range("A3").Value = "Naina"
$ws.Range("A4").Value = "Mohit"
$ws.Range("A5").Value = "Devanshu"
